$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix truncated / "fake data" text values in column B (NEName)
#    (the real content edits called out by the commit: "update fake data")
# ---------------------------------------------------------------------------
$ws.Range("B2").Value  = "DJKLERWOOBO R"
$ws.Range("B5").Value  = "DET Eiewo 2E - J"
$ws.Range("B12").Value = "DJEEWE GTEEG 4"
$ws.Range("B13").Value = "DJKLERWOOBO R"
$ws.Range("B17").Value = "UJ EL Gfdeer - O"
$ws.Range("B18").Value = "RE Or Egldfg H"
$ws.Range("B19").Value = "C052B_C00922"
$ws.Range("B26").Value = "YUOYIUGT"

# ---------------------------------------------------------------------------
# 2. Append four new data rows (28-31), copying the formatting from an
#    existing fully-populated data row (row 8) so the new cells pick up the
#    same styles used throughout the table.
# ---------------------------------------------------------------------------
$ws.Range("A8:D8").Copy()
$ws.Range("A28:D28").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A29:D29").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A30:D30").PasteSpecial(-4122)
$ws.Range("A8:D8").Copy()
$ws.Range("A31:D31").PasteSpecial(-4122)

$ws.Range("F8").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F31").PasteSpecial(-4122)

$ws.Range("A28").Value = "UVUUD_SLEEF"
$ws.Range("B28").Value = "GJKLJLEE"
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = "BSC"

$ws.Range("A29").Value = "UVUUD_SLEEF"
$ws.Range("B29").Value = "BKGIEIEW"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "nodeB"

$ws.Range("A30").Value = "GGEDX_BGREW"
$ws.Range("B30").Value = "PVOSWIRF"
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = "RNC"

$ws.Range("A31").Value = "GGEDX_BGREW"
$ws.Range("B31").Value = "VDHJWQO"
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = "nodeB"

# ---------------------------------------------------------------------------
# 3. Re-apply the VLOOKUP formula across F8:F31 so the newly extended range
#    is rebuilt as one shared formula (matches the xml diff's t="shared").
# ---------------------------------------------------------------------------
$ws.Range("F8:F31").Formula = "=VLOOKUP(E:E,B:B,1,0)"

# ---------------------------------------------------------------------------
# 4. Conditional formatting: split the old "duplicate values" rule so it no
#    longer covers B1, and add a fresh "duplicate values" rule scoped to B1
#    only (same red-on-red look), matching the dxfs/conditionalFormatting
#    changes in the diff.
# ---------------------------------------------------------------------------
$oldRule = $ws.Range("B1:B7").FormatConditions.Item(1)
$oldRule.ModifyAppliesToRange($ws.Range("B2:B7"))
$oldRule.Priority = 1

$newRule = $ws.Range("B1").FormatConditions.AddUniqueValues()
$newRule.DupeUnique = 1
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615
$newRule.Priority = 3

# ---------------------------------------------------------------------------
# 5. Restore the author's last selection.
# ---------------------------------------------------------------------------
$ws.Range("B13").Select()
